$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 227.7
$ws.Range("I11").Value = 227.7
$ws.Range("K11").Value = 227.7
$ws.Range("M11").Value = -87.69999999999999

$ws.Range("H40").Value = 2515
$ws.Range("I40").Value = 1460
$ws.Range("J40").Value = 3632.0588
$ws.Range("K40").Value = 1460
$ws.Range("L40").Value = 3632.0588
$ws.Range("M40").Value = -1285
$ws.Range("N40").Value = -3982.0588

$ws.Range("H70").Value = 2513.3333
$ws.Range("I70").Value = 1270
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 3810
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -3540
$ws.Range("N70").Value = -15540

$ws.Range("H73").Value = 2513.3333
$ws.Range("I73").Value = 1270
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 3810
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -2874
$ws.Range("N73").Value = -16872

$ws.Range("H125").Value = 3378.875
$ws.Range("I125").Value = 2516
$ws.Range("J125").Value = 3666.5
$ws.Range("K125").Value = 22644
$ws.Range("L125").Value = 32998.5
$ws.Range("M125").Value = -20184
$ws.Range("N125").Value = -37918.5

$ws.Range("H129").Value = 874
$ws.Range("J129").Value = 1041
$ws.Range("L129").Value = 3123
$ws.Range("N129").Value = -13123

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4306.288
$ws.Range("I32").Value = 2779.3137
$ws.Range("K32").Value = 2779.3137
$ws.Range("M32").Value = -2492.3137

$ws.Range("H61").Value = 1329.2593
$ws.Range("I61").Value = 1329.2593
$ws.Range("K61").Value = 1329.2593
$ws.Range("M61").Value = -1117.2593

$ws.Range("H125").Value = 35000
$ws.Range("J125").Value = 35000
$ws.Range("L125").Value = 35000
$ws.Range("N125").Value = -44840

$ws.Range("H128").Value = 34000
$ws.Range("J128").Value = 34000
$ws.Range("L128").Value = 34000
$ws.Range("N128").Value = -43960

$ws.Range("H131").Value = 33315.31
$ws.Range("J131").Value = 33315.31
$ws.Range("L131").Value = 33315.31
$ws.Range("N131").Value = -43395.31

$ws.Range("H132").Value = 1884.12
$ws.Range("I132").Value = 1590.6666
$ws.Range("J132").Value = 3424.75
$ws.Range("K132").Value = 4771.9998
$ws.Range("L132").Value = 10274.25
$ws.Range("M132").Value = -2241.9998
$ws.Range("N132").Value = -15334.25

$ws.Range("H136").Value = 1329.2593
$ws.Range("I136").Value = 1329.2593
$ws.Range("K136").Value = 3987.7779
$ws.Range("M136").Value = -1437.7779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 35000
$ws.Range("J126").Value = 35000
$ws.Range("L126").Value = 35000
$ws.Range("N126").Value = -44880

$ws.Range("H130").Value = 515000
$ws.Range("J130").Value = 515000
$ws.Range("L130").Value = 515000
$ws.Range("N130").Value = -525040

$ws.Range("H134").Value = 4643.6
$ws.Range("I134").Value = 4740.607
$ws.Range("J134").Value = 4417.25
$ws.Range("K134").Value = 14221.821
$ws.Range("L134").Value = 13251.75
$ws.Range("M134").Value = -11686.821
$ws.Range("N134").Value = -18321.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 38333.11
$ws.Range("J20").Value = 38333.11
$ws.Range("L20").Value = 38333.11
$ws.Range("N20").Value = -38805.11

$ws.Range("H30").Value = 38333.11
$ws.Range("J30").Value = 38333.11
$ws.Range("L30").Value = 38333.11
$ws.Range("N30").Value = -38515.11

$ws.Range("H60").Value = 13632
$ws.Range("I60").Value = 6900
$ws.Range("J60").Value = 16998
$ws.Range("K60").Value = 6900
$ws.Range("L60").Value = 16998
$ws.Range("M60").Value = -6389
$ws.Range("N60").Value = -18020

$ws.Range("H68").Value = 16500

$ws.Range("H71").Value = 16500

$ws.Range("H74").Value = 19657
$ws.Range("J74").Value = 19657
$ws.Range("L74").Value = 19657
$ws.Range("N74").Value = -21405

$ws.Range("H77").Value = 19657
$ws.Range("J77").Value = 19657
$ws.Range("L77").Value = 58971
$ws.Range("N77").Value = -67707

$ws.Range("H127").Value = 33799
$ws.Range("J127").Value = 33799
$ws.Range("L127").Value = 33799
$ws.Range("N127").Value = -43719

$ws.Range("H128").Value = 38333.11
$ws.Range("J128").Value = 38333.11
$ws.Range("L128").Value = 38333.11
$ws.Range("N128").Value = -48293.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2012.2858
$ws.Range("J21").Value = 1490
$ws.Range("L21").Value = 4470
$ws.Range("N21").Value = -4816

$ws.Range("H48").Value = 1993.2632
$ws.Range("J48").Value = 1993.2632
$ws.Range("L48").Value = 5979.7896
$ws.Range("N48").Value = -6479.7896

$ws.Range("H107").Value = 542384.9399999999
$ws.Range("I107").Value = 499.79166
$ws.Range("J107").Value = 774621.4
$ws.Range("K107").Value = 1499.37498
$ws.Range("L107").Value = 2323864.2
$ws.Range("M107").Value = 420.6250199999999
$ws.Range("N107").Value = -2327704.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 34649.75
$ws.Range("J124").Value = 34649.75
$ws.Range("L124").Value = 34649.75
$ws.Range("N124").Value = -44469.75

$ws.Range("H130").Value = 1000000
$ws.Range("J130").Value = 1000000
$ws.Range("L130").Value = 1000000
$ws.Range("N130").Value = -1010040

$ws.Range("H132").Value = 1847.4865
$ws.Range("I132").Value = 1701.5358
$ws.Range("J132").Value = 2301.5557
$ws.Range("K132").Value = 5104.607400000001
$ws.Range("L132").Value = 6904.6671
$ws.Range("M132").Value = -2574.607400000001
$ws.Range("N132").Value = -11964.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 7600.75
$ws.Range("J38").Value = 7600.75
$ws.Range("L38").Value = 7600.75
$ws.Range("N38").Value = -8420.75

$ws.Range("H127").Value = 56000
$ws.Range("J127").Value = 56000
$ws.Range("L127").Value = 56000
$ws.Range("N127").Value = -65920

$ws.Range("H130").Value = 29800
$ws.Range("J130").Value = 29800
$ws.Range("L130").Value = 29800
$ws.Range("N130").Value = -39840

$ws.Range("H132").Value = 4110.952
$ws.Range("I132").Value = 3957.5
$ws.Range("J132").Value = 4602
$ws.Range("K132").Value = 11872.5
$ws.Range("L132").Value = 13806
$ws.Range("M132").Value = -9342.5
$ws.Range("N132").Value = -18866

$ws.Range("H136").Value = 1511.2572
$ws.Range("I136").Value = 1415.2
$ws.Range("J136").Value = 1639.3334
$ws.Range("K136").Value = 4245.6
$ws.Range("L136").Value = 4918.0002
$ws.Range("M136").Value = -1695.6
$ws.Range("N136").Value = -10018.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 3633.8333
$ws.Range("I55").Value = 7000
$ws.Range("J55").Value = 2960.6
$ws.Range("K55").Value = 7000
$ws.Range("L55").Value = 2960.6
$ws.Range("M55").Value = -6723
$ws.Range("N55").Value = -3514.6

$ws.Range("H113").Value = 478.14285
$ws.Range("I113").Value = 256.2857
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 768.8571000000001
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1401.1429
$ws.Range("N113").Value = -6440
